$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly record: insert a row at 13, pushing the existing rows 13-39 down
# to 14-40 (dimension grows from A1:R39 to A1:R40).
$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value = 2
$ws.Cells.Item(13, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(13, 3).Value = "Coquimbo"
$ws.Cells.Item(13, 4).Value = "10/13/2021"
$ws.Cells.Item(13, 5).Value = 4
$ws.Cells.Item(13, 6).Value = 100112026
$ws.Cells.Item(13, 7).Value = "Haba"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 1600
$ws.Cells.Item(13, 11).Value = 4000
$ws.Cells.Item(13, 12).Value = 5000
$ws.Cells.Item(13, 13).Value = 4500
$ws.Cells.Item(13, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(13, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(13, 16).Value = 180
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = "Hortaliza"
